$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text: lesson_ID -> lesson_id
$ws.Range("A1").Value = "lesson_id"

# Shift lesson ids down by one (1..14 -> 0..13)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Update the active selection to D2:D15 with active cell D2
$ws.Range("D2:D15").Select()
